$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0003933774834437086
$ws.Range("C2").Value = 0.0003915181315304241
$ws.Range("D2").Value = 0.0003985354350254437
$ws.Range("E2").Value = 0.0004071883530482257
$ws.Range("F2").Value = 0.0004163454124903623
$ws.Range("G2").Value = 0.0004096987951807229
$ws.Range("H2").Value = 0.0004116953762466002
$ws.Range("B3").Value = 0.006193377483443708
$ws.Range("C3").Value = 0.006306084818684696
$ws.Range("D3").Value = 0.006228124612138513
$ws.Range("E3").Value = 0.006528662420382165
$ws.Range("F3").Value = 0.006626831148804934
$ws.Range("G3").Value = 0.006656626506024096
$ws.Range("H3").Value = 0.006533998186763373
$ws.Range("B4").Value = 0.09928476821192053
$ws.Range("C4").Value = 0.1011063306699447
$ws.Range("D4").Value = 0.1022961400024823
$ws.Range("E4").Value = 0.1044358507734304
$ws.Range("F4").Value = 0.1058596761757903
$ws.Range("G4").Value = 0.1031927710843373
$ws.Range("H4").Value = 0.1058930190389846
$ws.Range("B5").Value = 0.4998675496688742
$ws.Range("C5").Value = 0.5144437615242777
$ws.Range("D5").Value = 0.5150800546108973
$ws.Range("E5").Value = 0.5286624203821656
$ws.Range("F5").Value = 0.5366229760986893
$ws.Range("G5").Value = 0.5283734939759036
$ws.Range("H5").Value = 0.5332728921124207
$ws.Range("B7").Value = 1.587019867549669
$ws.Range("C7").Value = 1.606023355869699
$ws.Range("D7").Value = 1.645773861238674
$ws.Range("E7").Value = 1.669699727024568
$ws.Range("F7").Value = 1.69313801079414
$ws.Range("G7").Value = 1.676506024096386
$ws.Range("H7").Value = 1.694469628286491
$ws.Range("B8").Value = 3.856953642384106
$ws.Range("C8").Value = 3.884449907805777
$ws.Range("D8").Value = 4.00893632865831
$ws.Range("E8").Value = 4.076433121019108
$ws.Range("F8").Value = 4.117193523515806
$ws.Range("G8").Value = 4.119879518072289
$ws.Range("H8").Value = 4.115140525838622
$ws.Range("B9").Value = 8.055629139072847
$ws.Range("C9").Value = 8.168408113091578
$ws.Range("D9").Value = 8.321956063050763
$ws.Range("E9").Value = 8.47361237488626
$ws.Range("F9").Value = 8.600616808018504
$ws.Range("G9").Value = 8.445783132530121
$ws.Range("H9").Value = 8.556663644605621
$ws.Range("B10").Value = 14.25960264900662
$ws.Range("C10").Value = 15.09526736324523
$ws.Range("D10").Value = 15.41516693558396
$ws.Range("E10").Value = 15.7393084622384
$ws.Range("F10").Value = 15.7710100231303
$ws.Range("G10").Value = 15.7289156626506
$ws.Range("H10").Value = 15.90208522212149
